$d = $word.ActiveDocument

$pairs = @(
    ,@('2025-02-04 Tuesday', '2025-02-05 Wednesday')
    ,@('77+17=94', '95-46=49')
    ,@('18+4=22', '8+34=42')
    ,@('62-49=13', '42-29=13')
    ,@('16+27=43', '93-88=5')
    ,@('71-9=62', '95-48=47')
    ,@('8+8=16', '71-15=56')
    ,@('42-15=27', '76+18=94')
    ,@('49+22=71', '56+28=84')
    ,@('17+5=22', '94-25=69')
    ,@('82-38=44', '35-6=29')
    ,@('58+33=91', '8+17=25')
    ,@('64-26=38', '29+24=53')
    ,@('52+9=61', '71-29=42')
    ,@('49+49=98', '56-27=29')
    ,@('58+37=95', '29+19=48')
    ,@('76-48=28', '82-13=69')
    ,@('49+19=68', '45+49=94')
    ,@('85-78=7', '14+19=33')
    ,@('71-47=24', '83-36=47')
    ,@('37+56=93', '13+49=62')
    ,@('75+9=84', '27-18=9')
    ,@('35+47=82', '48+37=85')
    ,@('6+75=81', '92-44=48')
    ,@('49+14=63', '82-18=64')
    ,@('94-26=68', '55-39=16')
    ,@('92-53=39', '25+27=52')
    ,@('96-39=57', '67+29=96')
    ,@('8+56=64', '40-8=32')
    ,@('40-6=34', '92-48=44')
    ,@('73-39=34', '37-9=28')
    ,@('84-76=8', '38+8=46')
    ,@('47+8=55', '19+2=21')
    ,@('22+9=31', '87-49=38')
    ,@('80-46=34', '96-59=37')
    ,@('47+15=62', '90-21=69')
    ,@('68-59=9', '83-79=4')
    ,@('7+49=56', '72-56=16')
    ,@('70-44=26', '46+9=55')
    ,@('17+35=52', '42-35=7')
    ,@('68+4=72', '51-5=46')
    ,@('19+39=58', '17+18=35')
    ,@('28+49=77', '32-15=17')
    ,@('9+42=51', '60-39=21')
    ,@('74-37=37', '88-79=9')
    ,@('72-27=45', '88-29=59')
    ,@('15+46=61', '81-8=73')
    ,@('83+8=91', '72-69=3')
    ,@('22+49=71', '38+9=47')
    ,@('30-19=11', '53-18=35')
    ,@('73-26=47', '4+79=83')
    ,@('41-33=8', '49+39=88')
    ,@('8+63=71', '55+18=73')
    ,@('6+65=71', '51-28=23')
    ,@('39+48=87', '59+13=72')
    ,@('93-38=55', '19+28=47')
    ,@('39+37=76', '7+27=34')
    ,@('92-43=49', '71-69=2')
    ,@('78-69=9', '60-38=22')
    ,@('40-15=25', '49+45=94')
    ,@('7+66=73', '72-26=46')
    ,@('80-24=56', '23-15=8')
    ,@('84-26=58', '19+76=95')
    ,@('83-55=28', '44+19=63')
    ,@('36-7=29', '14+39=53')
    ,@('54+28=82', '15+39=54')
    ,@('26+48=74', '14+17=31')
    ,@('23-19=4', '24+69=93')
    ,@('80-36=44', '55+8=63')
    ,@('9+77=86', '37+7=44')
    ,@('71-26=45', '48+4=52')
    ,@('67-18=49', '47+18=65')
    ,@('16+25=41', '26+67=93')
    ,@('65+7=72', '14+49=63')
    ,@('64-5=59', '93-85=8')
    ,@('22-13=9', '91-68=23')
    ,@('90-41=49', '68+14=82')
    ,@('60-34=26', '35-27=8')
    ,@('39+5=44', '67+7=74')
    ,@('8+39=47', '39+34=73')
    ,@('57+36=93', '7+65=72')
    ,@('80-62=18', '31-22=9')
    ,@('26+9=35', '3+88=91')
    ,@('71-34=37', '7+34=41')
    ,@('51-29=22', '9+84=93')
    ,@('84+8=92', '8+68=76')
    ,@('94-46=48', '6+55=61')
    ,@('65-7=58', '26+27=53')
    ,@('97-29=68', '18+49=67')
    ,@('36+29=65', '77+4=81')
    ,@('82-4=78', '8+86=94')
    ,@('17+25=42', '8+53=61')
    ,@('17+28=45', '80-42=38')
    ,@('13+79=92', '61-12=49')
    ,@('28-9=19', '51-7=44')
    ,@('30-21=9', '80-35=45')
    ,@('50-24=26', '74+7=81')
    ,@('44+38=82', '66+19=85')
    ,@('42-8=34', '9+76=85')
    ,@('71-13=58', '94-5=89')
    ,@('55-9=46', '14+39=53')
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

Write-Host "Done"
